$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '26.412.71'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '1.796.59'
$ws.Range("E3").Value = '  -1.95%  '
Set-TextValue "D4" '1.009'
$ws.Range("E4").Value = '  +0.24%  '
Set-TextValue "D5" '1.008'
$ws.Range("E5").Value = '  +0.21%  '
Set-TextValue "D6" '308.01'
$ws.Range("E6").Value = '  -0.83%  '
Set-TextValue "D7" '0.4525'
$ws.Range("E7").Value = '  -2.00%  '
Set-TextValue "D8" '0.3595'
$ws.Range("E8").Value = '  -2.04%  '
Set-TextValue "D9" '46.37'
$ws.Range("E9").Value = '  +1.15%  '
Set-TextValue "D10" '0.07124'
$ws.Range("E10").Value = '  -0.71%  '
Set-TextValue "D11" '0.8881'
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("E12").Value = '  -0.50%  '
Set-TextValue "D13" '19.53'
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").Value = '1.794.86'
$ws.Range("E14").Value = '  -2.21%  '
Set-TextValue "D15" '5.288'
$ws.Range("E15").Value = '  -0.88%  '
Set-TextValue "D16" '6.337'
$ws.Range("E16").Value = '  -0.84%  '
Set-TextValue "D17" '85.04'
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("E18").Value = '  +0.27%  '
Set-TextValue "D19" '0.000008579'
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  -1.15%  '
$ws.Range("D22").Value = '26.440.74'
$ws.Range("E22").Value = '  -1.78%  '
Set-TextValue "D23" '4.998'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.057.68'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D25" '10.53'
$ws.Range("E25").Value = '  +0.76%  '
Set-TextValue "D26" '1.985'
$ws.Range("E26").Value = '  +0.50%  '
Set-TextValue "D27" '152.88'
$ws.Range("E27").Value = '  +1.43%  '
Set-TextValue "D28" '17.93'
$ws.Range("E28").Value = '  -1.66%  '
Set-TextValue "D29" '2.043'
$ws.Range("E29").Value = '  +3.61%  '
Set-TextValue "D30" '112.19'
$ws.Range("E30").Value = '  -1.23%  '
Set-TextValue "D31" '4.869'
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("E32").Value = '  -1.79%  '
Set-TextValue "D33" '3.049'
$ws.Range("E33").Value = '  -2.57%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D34" '4.460'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D35" '0.7300'
$ws.Range("E35").Value = '  -3.52%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D36" '2.724'
$ws.Range("E36").Value = '  +5.28%  '
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("E38").Value = '  -1.64%  '
Set-TextValue "D39" '0.01936'
$ws.Range("E39").Value = '  +0.11%  '
Set-TextValue "D40" '0.05108'
$ws.Range("E40").Value = '  -0.53%  '
Set-TextValue "D41" '2.879'
$ws.Range("E41").Value = '  -1.70%  '
Set-TextValue "D42" '0.5164'
$ws.Range("E42").Value = '  +3.66%  '
Set-TextValue "D43" '6.898'
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("E44").Value = '  -4.53%  '
Set-TextValue "D45" '8.003'
$ws.Range("E45").Value = '  -4.26%  '
Set-TextValue "D46" '0.4672'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("E47").Value = '  +0.21%  '
Set-TextValue "D48" '9.876'
$ws.Range("E48").Value = '  -2.76%  '
Set-TextValue "D49" '100.94'
$ws.Range("E49").Value = '  -1.35%  '
$ws.Range("E50").Value = '  -1.48%  '
Set-TextValue "D51" '64.56'
$ws.Range("E51").Value = '  +0.15%  '
